# Change clients import xlsx file format
# - Replace the header row of the single worksheet with the new, smaller set
#   of columns (category_client, account_type, account_number, firstname,
#   lastname, sexe, telephone, email, ville)
# - Drop the old extra rows / cell formatting (bold headers, number format,
#   custom column widths, page setup) that belonged to the previous layout
# - Update the selected cell in the sheet view
# - Swap the accent1 / accent5 theme colors

$wb = $excel.ActiveWorkbook

# Start from a brand-new, completely blank worksheet so none of the old
# column widths / cell styles / extra rows survive, then put it in the
# place of (and rename it to) the original sheet.
$oldSheetName = $wb.ActiveSheet.Name
$newSheet = $wb.Worksheets.Add()
$oldSheet = $wb.Worksheets.Item($oldSheetName)
$oldSheet.Delete() | Out-Null
$newSheet.Name = $oldSheetName

# New header row
$headers = @(
    "category_client",
    "account_type",
    "account_number",
    "firstname",
    "lastname",
    "sexe",
    "telephone",
    "email",
    "ville"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the saved selection/cursor position
$newSheet.Range("C6").Select() | Out-Null

# Swap theme accent1 <-> accent5 colors (5B9BD5 <-> 4472C4)
$scheme = $wb.Theme.ThemeColorScheme
$accent1 = $scheme.Colors(5)
$accent5 = $scheme.Colors(9)
$old1 = $accent1.RGB
$old5 = $accent5.RGB
$accent1.RGB = $old5
$accent5.RGB = $old1
